$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price (column D) and Volume(1h) (column E) updates
# from the Sun Jan 15 06:33:27 UTC 2023 GitHub Actions symbol-list refresh.
# A leading apostrophe forces text entry so the numeric-looking strings
# (and any significant trailing zeros) are preserved exactly as text,
# matching the inline-string cell type already used throughout the sheet.

$ws.Cells.Item(2, 4).Value = "'298.22"
$ws.Cells.Item(2, 5).Value = "'-3.32%"

$ws.Cells.Item(3, 4).Value = "'31.94"
$ws.Cells.Item(3, 5).Value = "'-0.71%"

$ws.Cells.Item(4, 4).Value = "'5.114"
$ws.Cells.Item(4, 5).Value = "'-4.07%"

$ws.Cells.Item(5, 4).Value = "'0.07541"
$ws.Cells.Item(5, 5).Value = "'0.87%"

$ws.Cells.Item(6, 4).Value = "'7.754"
$ws.Cells.Item(6, 5).Value = "'-0.63%"

$ws.Cells.Item(7, 4).Value = "'1.746"
$ws.Cells.Item(7, 5).Value = "'10.79%"

$ws.Cells.Item(8, 4).Value = "'3.793"
$ws.Cells.Item(8, 5).Value = "'3.40%"

$ws.Cells.Item(9, 4).Value = "'0.9293"
$ws.Cells.Item(9, 5).Value = "'2.46%"

$ws.Cells.Item(10, 5).Value = "'1.08%"

$ws.Cells.Item(11, 4).Value = "'0.07486"
$ws.Cells.Item(11, 5).Value = "'-2.42%"

$ws.Cells.Item(12, 4).Value = "'0.07929"
$ws.Cells.Item(12, 5).Value = "'-1.96%"

$ws.Cells.Item(13, 4).Value = "'0.03056"
$ws.Cells.Item(13, 5).Value = "'0.64%"

$ws.Cells.Item(14, 4).Value = "'0.09891"
$ws.Cells.Item(14, 5).Value = "'0.30%"

$ws.Cells.Item(15, 4).Value = "'0.001495"
$ws.Cells.Item(15, 5).Value = "'-1.43%"

$ws.Cells.Item(16, 4).Value = "'0.006489"
$ws.Cells.Item(16, 5).Value = "'1.20%"

$ws.Cells.Item(17, 4).Value = "'3.461"
$ws.Cells.Item(17, 5).Value = "'-0.79%"

$ws.Cells.Item(18, 4).Value = "'2.220"
$ws.Cells.Item(18, 5).Value = "'-0.88%"

$ws.Cells.Item(19, 4).Value = "'0.3281"
$ws.Cells.Item(19, 5).Value = "'0.42%"

$ws.Cells.Item(20, 5).Value = "'-0.73%"

$ws.Cells.Item(21, 4).Value = "'4.551"
$ws.Cells.Item(21, 5).Value = "'8.91%"

$ws.Cells.Item(22, 4).Value = "'0.04650"
$ws.Cells.Item(22, 5).Value = "'2.34%"

$ws.Cells.Item(23, 4).Value = "'0.1557"
$ws.Cells.Item(23, 5).Value = "'-3.88%"

$ws.Cells.Item(24, 4).Value = "'0.001221"
$ws.Cells.Item(24, 5).Value = "'0.45%"

$ws.Cells.Item(25, 4).Value = "'0.004423"
$ws.Cells.Item(25, 5).Value = "'-1.80%"

$ws.Cells.Item(27, 5).Value = "'6.79%"

$ws.Cells.Item(39, 4).Value = "'0.01680"
$ws.Cells.Item(39, 5).Value = "'-1.55%"

$ws.Cells.Item(40, 4).Value = "'0.04536"
$ws.Cells.Item(40, 5).Value = "'0.07%"

$ws.Cells.Item(41, 4).Value = "'0.007052"
$ws.Cells.Item(41, 5).Value = "'-1.61%"

$ws.Cells.Item(42, 4).Value = "'0.1326"
$ws.Cells.Item(42, 5).Value = "'-2.67%"

$ws.Cells.Item(43, 4).Value = "'0.002058"
$ws.Cells.Item(43, 5).Value = "'-8.80%"

$ws.Cells.Item(44, 4).Value = "'0.01166"
$ws.Cells.Item(44, 5).Value = "'-16.33%"

$ws.Cells.Item(45, 4).Value = "'0.00005973"
$ws.Cells.Item(45, 5).Value = "'-2.41%"

$ws.Cells.Item(46, 5).Value = "'1.34%"

$ws.Cells.Item(47, 5).Value = "'-0.19%"
